$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above row 216, shifting the existing 216:235 block down to 218:237
$ws.Range("A216:A217").EntireRow.Insert()

# Populate the two newly inserted rows with the new weekly price-report entries
$ws.Cells.Item(216, 1).Value = 11
$ws.Cells.Item(216, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(216, 3).Value = "Bíobío"
$ws.Cells.Item(216, 4).Value = 44449
$ws.Cells.Item(216, 5).Value = 8
$ws.Cells.Item(216, 6).Value = 100112004
$ws.Cells.Item(216, 7).Value = "Cebolla"
$ws.Cells.Item(216, 8).Value = "Sin especificar"
$ws.Cells.Item(216, 9).Value = "1a (guarda)"
$ws.Cells.Item(216, 10).Value = 700
$ws.Cells.Item(216, 11).Value = 6000
$ws.Cells.Item(216, 12).Value = 6500
$ws.Cells.Item(216, 13).Value = 6286
$ws.Cells.Item(216, 14).Value = "$/malla 18 kilos"
$ws.Cells.Item(216, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(216, 16).Value = 349
$ws.Cells.Item(216, 17).Value = 18
$ws.Cells.Item(216, 18).Value = "Hortaliza"

$ws.Cells.Item(217, 1).Value = 11
$ws.Cells.Item(217, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(217, 3).Value = "Bíobío"
$ws.Cells.Item(217, 4).Value = 44449
$ws.Cells.Item(217, 5).Value = 8
$ws.Cells.Item(217, 6).Value = 100112004
$ws.Cells.Item(217, 7).Value = "Cebolla"
$ws.Cells.Item(217, 8).Value = "Sin especificar"
$ws.Cells.Item(217, 9).Value = "2a (guarda)"
$ws.Cells.Item(217, 10).Value = 300
$ws.Cells.Item(217, 11).Value = 5000
$ws.Cells.Item(217, 12).Value = 5000
$ws.Cells.Item(217, 13).Value = 5000
$ws.Cells.Item(217, 14).Value = "$/malla 18 kilos"
$ws.Cells.Item(217, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(217, 16).Value = 278
$ws.Cells.Item(217, 17).Value = 18
$ws.Cells.Item(217, 18).Value = "Hortaliza"
